# Add a second 'Logging' timesheet to the workbook.
#
# Before: one sheet "Tabelle1" holding the timesheet grid.
# After : two sheets -
#   1) "Logging"   - small carryover bookkeeping sheet (new, first/leftmost)
#   2) "Timesheet" - the original "Tabelle1" grid, renamed, still active.

$wb = $excel.ActiveWorkbook

# A freshly added sheet is inserted right before the (currently only)
# sheet, which puts it first - exactly where "Logging" belongs. Do this
# BEFORE grabbing/renaming the original sheet: sheet handles here resolve
# by live position, and Add() shifts every existing sheet one slot over.
$logging = $wb.Worksheets.Add()
$logging.Name = "Logging"

# Now fetch the pre-existing sheet by its (still valid) original name and
# rename it.
$timesheet = $wb.Worksheets.Item("Tabelle1")
$timesheet.Name = "Timesheet"

# Small carryover log: row/column pointers used by the app to know where
# the last entry was written, plus the hours carried over.
$logging.Range("B1").Value = "carryover"
$logging.Range("A2").Value = "row"
$logging.Range("B2").Value = 33
$logging.Range("A3").Value = "column"
$logging.Range("B3").Value = 10

# The Timesheet tab stays the active / selected one, with the selection
# collapsed down to the single cell D12 (was D12:G31).
$timesheet.Activate() | Out-Null
$timesheet.Range("D12").Select() | Out-Null
